$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells
$ws.Range("B3").Value = 15
$ws.Range("B4").Value = 7

# Add new rows of data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 200

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 90

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 0

# Update selection to reflect the next empty cell
$ws.Range("B9").Select()
